$wb = $excel.ActiveWorkbook

# --- NuevaHoja (first sheet): replace its old "A22 = Hoja 1" content with the
# new lesson table: A1 = wrapped title text, B1 = integer, C1 = long-format date,
# A2 = a simple formula. ---
$ws1 = $wb.Worksheets.Item(1)

# Drop the old content in A22 ("Hoja 1") entirely - it disappears from the sheet.
$ws1.Range("A22").Value = $null

# A1: two-line wrapped text.
$ws1.Range("A1").Value = "Primer `nTexto"
$ws1.Range("A1").WrapText = $true

# B1: a plain integer, vertically centered.
$ws1.Range("B1").Value = 1000
$ws1.Range("B1").NumberFormat = "0"
$ws1.Range("B1").VerticalAlignment = -4108

# C1: a date (2022-01-01, serial 44562) shown in the Windows long-date format,
# underlined like the workbook's secondary font.
$ws1.Range("C1").Value = 44562
$ws1.Range("C1").NumberFormat = "[$-F800]dddd, mmmm dd, yyyy"
$ws1.Range("C1").Font.Underline = $true

# Column C is widened to fit the long date text.
$ws1.Range("C1").EntireColumn.ColumnWidth = 20.94

# A2: a tiny arithmetic formula.
$ws1.Range("A2").Formula = "=2+3"

# Scroll/select so C1 is the active cell, top-left of the view.
$ws1.Range("C1").Select()
$ws1.Application.ActiveWindow.ScrollColumn = 3
$ws1.Application.ActiveWindow.ScrollRow = 1

# Make NuevaHoja the active sheet/tab.
$ws1.Activate()

Write-Output "done"
